# Updated symbol list on Fri Feb 17 03:39:01 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Rows 2-16: price/volume updates only
Set-TextCell "D2" "309.46"
Set-TextCell "E2" "-3.90%"

Set-TextCell "D3" "48.41"
Set-TextCell "E3" "-7.55%"

Set-TextCell "D4" "5.155"
Set-TextCell "E4" "-4.10%"

Set-TextCell "D5" "0.07769"
Set-TextCell "E5" "-4.11%"

Set-TextCell "D6" "4.475"
Set-TextCell "E6" "-2.20%"

Set-TextCell "D7" "1.314"
Set-TextCell "E7" "18.83%"

Set-TextCell "D8" "1.560"
Set-TextCell "E8" "-5.43%"

Set-TextCell "E9" "-6.81%"

Set-TextCell "D10" "0.1941"
Set-TextCell "E10" "-0.74%"

Set-TextCell "D11" "0.04682"
Set-TextCell "E11" "2.92%"

Set-TextCell "D12" "0.09284"
Set-TextCell "E12" "-3.70%"

Set-TextCell "E13" "0.10%"

Set-TextCell "D14" "0.001260"
Set-TextCell "E14" "-4.92%"

Set-TextCell "D15" "0.04173"
Set-TextCell "E15" "-2.95%"

Set-TextCell "D16" "0.005814"
Set-TextCell "E16" "-0.62%"

# Rows 17-24: the coin list was re-ranked/shifted, names/links/prices changed
Set-TextCell "B17" "HotbitToken"
Set-TextCell "C17" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell "D17" "0.004125"
Set-TextCell "E17" "-3.33%"

Set-TextCell "B18" "LEO"
Set-TextCell "C18" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D18" "3.329"
Set-TextCell "E18" "-1.50%"

Set-TextCell "B19" "BTSEToken"
Set-TextCell "C19" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D19" "2.276"
Set-TextCell "E19" "-6.21%"

Set-TextCell "B20" "BitpandaEcosystemToken"
Set-TextCell "C20" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell "D20" "0.3492"
Set-TextCell "E20" "2.88%"

Set-TextCell "B21" "MCDex"
Set-TextCell "C21" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D21" "8.246"
Set-TextCell "E21" "0.71%"

Set-TextCell "B22" "ProBitToken"
Set-TextCell "C22" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell "D22" "0.1354"
Set-TextCell "E22" "-3.19%"

Set-TextCell "B23" "ZBToken"
Set-TextCell "C23" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell "D23" "0.3030"
Set-TextCell "E23" "3.52%"

Set-TextCell "B24" "BitKan"
Set-TextCell "C24" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell "D24" "0.001279"
Set-TextCell "E24" "-2.07%"

Set-TextCell "E25" "0.20%"
Set-TextCell "E26" "-3.89%"

Set-TextCell "E38" "-6.95%"

Set-TextCell "D39" "0.05877"
Set-TextCell "E39" "6.04%"

Set-TextCell "E40" "71.13%"

Set-TextCell "D41" "0.007934"
Set-TextCell "E41" "2.16%"

Set-TextCell "D42" "0.1423"
Set-TextCell "E42" "-1.51%"

Set-TextCell "D43" "0.008407"
Set-TextCell "E43" "9.45%"

Set-TextCell "D44" "0.007659"
Set-TextCell "E44" "-13.19%"

Set-TextCell "D45" "0.3119"
Set-TextCell "E45" "-11.48%"

Set-TextCell "D46" "0.00006958"
Set-TextCell "E46" "1.89%"

Set-TextCell "E47" "0.20%"

Set-TextCell "D48" "0.05667"
Set-TextCell "E48" "-6.47%"

Set-TextCell "E49" "0.31%"

Set-TextCell "D50" "0.00002101"
Set-TextCell "E50" "0.20%"

Set-TextCell "D51" "0.0002001"
Set-TextCell "E51" "0.20%"
